$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price column cells so values like "1.00" or
# "143.23" are preserved as literal text instead of being parsed as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values
$ws.Range("D2").Value = '64.439.12'
$ws.Range("E2").Value = '  +5.45%  '
$ws.Range("D3").Value = '3.080.90'
$ws.Range("E3").Value = '  +3.52%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("E5").Value = '  +2.00%  '
$ws.Range("D6").Value = '143.23'
$ws.Range("E6").Value = '  +10.00%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.074.03'
$ws.Range("E8").Value = '  +3.49%  '
$ws.Range("E9").Value = '  +0.89%  '
$ws.Range("D10").Value = '6.95'
$ws.Range("E10").Value = '  +18.21%  '
$ws.Range("E11").Value = '  +6.21%  '
$ws.Range("D12").Value = '0.455'
$ws.Range("E12").Value = '  +3.43%  '
$ws.Range("E13").Value = '  +4.29%  '
$ws.Range("D14").Value = '35.23'
$ws.Range("E14").Value = '  +5.66%  '
$ws.Range("D15").Value = '3.576.84'
$ws.Range("E15").Value = '  +3.41%  '
$ws.Range("D16").Value = '64.459.10'
$ws.Range("E16").Value = '  +5.50%  '
$ws.Range("D17").Value = '3.084.56'
$ws.Range("E17").Value = '  +3.81%  '
$ws.Range("E18").Value = '  -1.05%  '
$ws.Range("E19").Value = '  +3.09%  '
$ws.Range("D20").Value = '479.91'
$ws.Range("E20").Value = '  +2.16%  '
$ws.Range("D21").Value = '13.60'
$ws.Range("E21").Value = '  +5.07%  '
$ws.Range("D22").Value = '0.670'
$ws.Range("D23").Value = '7.54'
$ws.Range("E23").Value = '  +9.59%  '
$ws.Range("E24").Value = '  +11.06%  '
$ws.Range("D25").Value = '80.72'
$ws.Range("E25").Value = '  +1.88%  '
$ws.Range("E26").Value = '  +1.31%  '
$ws.Range("D27").Value = '2.80'
$ws.Range("E27").Value = '  +4.60%  '
$ws.Range("D28").Value = '7.89'
$ws.Range("E28").Value = '  +4.95%  '
$ws.Range("D29").Value = '2.05'
$ws.Range("E29").Value = '  +10.32%  '
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.16%  '
$ws.Range("D31").Value = '26.09'
$ws.Range("E31").Value = '  +3.36%  '
$ws.Range("E32").Value = '  +3.34%  '
$ws.Range("E33").Value = '  +7.18%  '
$ws.Range("D34").Value = '5.70'
$ws.Range("E34").Value = '  +5.71%  '
$ws.Range("D35").Value = '55.20'
$ws.Range("E35").Value = '  +1.62%  '
$ws.Range("D36").Value = '6.06'
$ws.Range("E36").Value = '  +4.70%  '
$ws.Range("D37").Value = '464.82'
$ws.Range("E37").Value = '  +5.12%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.0407'
$ws.Range("E38").Value = '  +9.03%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '0.0826'
$ws.Range("E39").Value = '  +5.75%  '
$ws.Range("D40").Value = '3.018.29'
$ws.Range("E40").Value = '  -2.93%  '
$ws.Range("D41").Value = '0.118'
$ws.Range("E41").Value = '  +2.13%  '
$ws.Range("D42").Value = '8.27'
$ws.Range("E42").Value = '  +3.63%  '
$ws.Range("D43").Value = '2.69'
$ws.Range("E43").Value = '  +18.83%  '
$ws.Range("D44").Value = '27.85'
$ws.Range("E44").Value = '  +11.11%  '
$ws.Range("E45").Value = '  +8.24%  '
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("E47").Value = '  +8.03%  '
$ws.Range("E48").Value = '  +4.31%  '
$ws.Range("E49").Value = '  +8.56%  '
$ws.Range("D50").Value = '117.15'
$ws.Range("E50").Value = '  +3.41%  '
$ws.Range("D51").Value = '2.06'
$ws.Range("E51").Value = '  +5.16%  '
